# Append one new data row (row 61) to the "Optical_Power" sheet, mirroring
# the structure of the existing rows (1..60): text columns A-L are stored as
# plain text (even when the text looks numeric, e.g. "-493", "15", "1"), and
# the coordinate columns M/N are stored as real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 61

# Helper-free approach: for cells whose text content could be mistaken for a
# number or date by Excel's input parser, force a text number format before
# assignment, then restore the cell style to "Normal" so no stray style index
# is left on the cell (matches the unstyled cells used by the source data).
function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "A$row" "-493"
Set-TextValue "B$row" "6/27/2025"
$ws.Range("C$row").Value = "JUFRE 424"
Set-TextValue "D$row" "15"
Set-TextValue "E$row" "807817955"
$ws.Range("F$row").Value = "Optical Power"
$ws.Range("G$row").Value = "Pendiente"
$ws.Range("H$row").Value = "Desmontar columna de 168 mm y traspasar redes a columna comunitaria"
Set-TextValue "I$row" "1"
$ws.Range("J$row").Value = "Desmonte"
$ws.Range("K$row").Value = "Sin equipos"
$ws.Range("L$row").Value = "Pasante"

$ws.Range("M$row").Value = -58.432644
$ws.Range("N$row").Value = -34.595434
